$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3489.2222
$ws.Range("I64").Value = 3566.6667
$ws.Range("J64").Value = 3450.5
$ws.Range("K64").Value = 3566.6667
$ws.Range("L64").Value = 3450.5
$ws.Range("M64").Value = -3318.6667
$ws.Range("N64").Value = -3946.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3489.2222
$ws.Range("I67").Value = 3566.6667
$ws.Range("J67").Value = 3450.5
$ws.Range("K67").Value = 3566.6667
$ws.Range("L67").Value = 3450.5
$ws.Range("M67").Value = -2708.6667
$ws.Range("N67").Value = -5166.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3916.611
$ws.Range("I74").Value = 4219.8
$ws.Range("J74").Value = 3800
$ws.Range("K74").Value = 4219.8
$ws.Range("L74").Value = 3800
$ws.Range("M74").Value = -3283.8
$ws.Range("N74").Value = -5672

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6131.773
$ws.Range("I76").Value = 4994.0586
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 4994.0586
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -4679.0586
$ws.Range("N76").Value = -10630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3916.611
$ws.Range("I77").Value = 4219.8
$ws.Range("J77").Value = 3800
$ws.Range("K77").Value = 21099
$ws.Range("L77").Value = 19000
$ws.Range("M77").Value = -16419
$ws.Range("N77").Value = -28360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6131.773
$ws.Range("I79").Value = 4994.0586
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 4994.0586
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -3902.0586
$ws.Range("N79").Value = -12184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 20000980
$ws.Range("I137").Value = 1000.8
$ws.Range("J137").Value = 40000960
$ws.Range("K137").Value = 3002.4
$ws.Range("L137").Value = 120002880
$ws.Range("M137").Value = -452.3999999999996
$ws.Range("N137").Value = -120007980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1079.5676
$ws.Range("I141").Value = 943.97144
$ws.Range("J141").Value = 3452.5
$ws.Range("K141").Value = 2831.91432
$ws.Range("L141").Value = 10357.5
$ws.Range("M141").Value = 2348.08568
$ws.Range("N141").Value = -20717.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 23813440
$ws.Range("I74").Value = 35716228
$ws.Range("J74").Value = 7861.143
$ws.Range("K74").Value = 35716228
$ws.Range("L74").Value = 7861.143
$ws.Range("M74").Value = -35715354
$ws.Range("N74").Value = -9609.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 23813440
$ws.Range("I77").Value = 35716228
$ws.Range("J77").Value = 7861.143
$ws.Range("K77").Value = 178581140
$ws.Range("L77").Value = 39305.715
$ws.Range("M77").Value = -178576772
$ws.Range("N77").Value = -48041.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2096.2856
$ws.Range("I102").Value = 1816.4445
$ws.Range("J102").Value = 2600
$ws.Range("K102").Value = 1816.4445
$ws.Range("L102").Value = 2600
$ws.Range("M102").Value = -194.4445000000001
$ws.Range("N102").Value = -5844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6099398
$ws.Range("I132").Value = 8930049
$ws.Range("J132").Value = 2610.923
$ws.Range("K132").Value = 26790147
$ws.Range("L132").Value = 7832.768999999999
$ws.Range("M132").Value = -26787617
$ws.Range("N132").Value = -12892.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1655.4445
$ws.Range("I20").Value = 1608.9524
$ws.Range("J20").Value = 1818.1666
$ws.Range("K20").Value = 1608.9524
$ws.Range("L20").Value = 1818.1666
$ws.Range("M20").Value = -1361.9524
$ws.Range("N20").Value = -2312.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5000.4546
$ws.Range("I105").Value = 4549.75
$ws.Range("J105").Value = 5100.6113
$ws.Range("K105").Value = 4549.75
$ws.Range("L105").Value = 5100.6113
$ws.Range("M105").Value = -2802.75
$ws.Range("N105").Value = -8594.6113

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4114.3
$ws.Range("I134").Value = 2527.6428
$ws.Range("J134").Value = 7816.5
$ws.Range("K134").Value = 7582.928400000001
$ws.Range("L134").Value = 23449.5
$ws.Range("M134").Value = -5047.928400000001
$ws.Range("N134").Value = -28519.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6670405.5
$ws.Range("I31").Value = 3767.0667
$ws.Range("J31").Value = 66670150
$ws.Range("K31").Value = 3767.0667
$ws.Range("L31").Value = 66670150
$ws.Range("M31").Value = -3472.0667
$ws.Range("N31").Value = -66670740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6670405.5
$ws.Range("I34").Value = 3767.0667
$ws.Range("J34").Value = 66670150
$ws.Range("K34").Value = 3767.0667
$ws.Range("L34").Value = 66670150
$ws.Range("M34").Value = -3565.0667
$ws.Range("N34").Value = -66670554

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14287774
$ws.Range("I132").Value = 16668209
$ws.Range("J132").Value = 5162.4
$ws.Range("K132").Value = 50004627
$ws.Range("L132").Value = 15487.2
$ws.Range("M132").Value = -50002097
$ws.Range("N132").Value = -20547.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 769514.25
$ws.Range("I134").Value = 1709.9445
$ws.Range("J134").Value = 1832627.9
$ws.Range("K134").Value = 5129.833500000001
$ws.Range("L134").Value = 5497883.699999999
$ws.Range("M134").Value = -2594.833500000001
$ws.Range("N134").Value = -5502953.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 9774.75
$ws.Range("I63").Value = 9699
$ws.Range("J63").Value = 9800
$ws.Range("K63").Value = 29097
$ws.Range("L63").Value = 29400
$ws.Range("M63").Value = -28348
$ws.Range("N63").Value = -30898

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 9774.75
$ws.Range("I66").Value = 9699
$ws.Range("J66").Value = 9800
$ws.Range("K66").Value = 87291
$ws.Range("L66").Value = 88200
$ws.Range("M66").Value = -83547
$ws.Range("N66").Value = -95688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11497812
$ws.Range("I80").Value = 16669432
$ws.Range("J80").Value = 5322.1113
$ws.Range("K80").Value = 16669432
$ws.Range("L80").Value = 5322.1113
$ws.Range("M80").Value = -16668434
$ws.Range("N80").Value = -7318.1113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 11497812
$ws.Range("I83").Value = 16669432
$ws.Range("J83").Value = 5322.1113
$ws.Range("K83").Value = 83347160
$ws.Range("L83").Value = 26610.5565
$ws.Range("M83").Value = -83342168
$ws.Range("N83").Value = -36594.5565

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 8212
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 8212
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 8212
$ws.Range("N104").Value = -15200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 18525504
$ws.Range("I136").Value = 26318574
$ws.Range("J136").Value = 16963
$ws.Range("K136").Value = 78955722
$ws.Range("L136").Value = 50889
$ws.Range("M136").Value = -78953172
$ws.Range("N136").Value = -55989

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3751.0454
$ws.Range("I132").Value = 2935
$ws.Range("J132").Value = 4730.3
$ws.Range("K132").Value = 8805
$ws.Range("L132").Value = 14190.9
$ws.Range("M132").Value = -6275
$ws.Range("N132").Value = -19250.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1476.4286
$ws.Range("I136").Value = 1186.1538
$ws.Range("J136").Value = 5250
$ws.Range("K136").Value = 3558.4614
$ws.Range("L136").Value = 15750
$ws.Range("M136").Value = -1008.4614
$ws.Range("N136").Value = -20850
